# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to match the regenerated gh-pages output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 570
$ws1.Range("F5").Value  = 1792
$ws1.Range("F6").Value  = 55
$ws1.Range("F8").Value  = 157
$ws1.Range("F9").Value  = 2150
$ws1.Range("F10").Value = 45
$ws1.Range("F11").Value = 150
$ws1.Range("F12").Value = 1362
$ws1.Range("F13").Value = 477
$ws1.Range("F14").Value = 26
$ws1.Range("F15").Value = 298
$ws1.Range("F16").Value = 212
$ws1.Range("F21").Value = 58
$ws1.Range("F22").Value = 19
$ws1.Range("F23").Value = 1169
$ws1.Range("F25").Value = 349
$ws1.Range("F27").Value = 275
$ws1.Range("F28").Value = 344

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 570
$ws4.Range("F5").Value  = 1792
$ws4.Range("F7").Value  = 55
$ws4.Range("F9").Value  = 157
$ws4.Range("F10").Value = 2150
$ws4.Range("F11").Value = 45
$ws4.Range("F12").Value = 150
$ws4.Range("F13").Value = 1362
$ws4.Range("F14").Value = 477
$ws4.Range("F15").Value = 26
$ws4.Range("F16").Value = 298
$ws4.Range("F17").Value = 212
$ws4.Range("F22").Value = 58
$ws4.Range("F23").Value = 19
$ws4.Range("F24").Value = 1169
$ws4.Range("F26").Value = 349
$ws4.Range("F28").Value = 275
$ws4.Range("F29").Value = 344
